$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (column C) for all existing data rows (2..408) from 45188 -> 45189
$ws.Range("C2:C408").Value = 45189

# 2. Row 408 picks up an explicit row height (ht="15" customHeight="1") in the diff.
$ws.Rows.Item(408).RowHeight = 15

# 3. Append the new row 409 with its data.
$ws.Cells.Item(409, 1).Value = "A 43904-2023"

$ws.Cells.Item(409, 2).Value = 45187
$ws.Cells.Item(409, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(409, 3).Value = 45189
$ws.Cells.Item(409, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(409, 4).Value = "ÖSTERGÖTLANDS LÄN"
$ws.Cells.Item(409, 5).Value = "YDRE"

$ws.Cells.Item(409, 7).Value = 4.5
$ws.Cells.Item(409, 8).Value = 0
$ws.Cells.Item(409, 9).Value = 0
$ws.Cells.Item(409, 10).Value = 0
$ws.Cells.Item(409, 11).Value = 0
$ws.Cells.Item(409, 12).Value = 0
$ws.Cells.Item(409, 13).Value = 0
$ws.Cells.Item(409, 14).Value = 0
$ws.Cells.Item(409, 15).Value = 0
$ws.Cells.Item(409, 16).Value = 0
$ws.Cells.Item(409, 17).Value = 0

$ws.Cells.Item(409, 18).Value = ""
$ws.Cells.Item(409, 18).WrapText = $true
